# CIERRE 7 ENE 2022
# - Update the "VALES DE INSENTIVOS" sheet's concept text from NOVIEMBRE to
#   DICIEMBRE 2021 (the amount text on ARQUITECTO stays the same wording).
# - Switch the active/selected tab from "ARQUITECTO" to "VALES DE INSENTIVOS",
#   and update that sheet's selected cell.

$wb = $excel.ActiveWorkbook

$wsVales = $wb.Worksheets.Item(2)

# Update the incentive-month wording on the "VALES DE INSENTIVOS" sheet.
$wsVales.Range("A4").Value = "PAGO DE INCENTIVO DEL MES DE DICIEMBRE  2021"

# Move the active tab / selection to "VALES DE INSENTIVOS", leaving
# "ARQUITECTO" unselected (its own in-sheet selection stays on G8).
$wsVales.Activate() | Out-Null
$wsVales.Range("H9").Select() | Out-Null
